# Auto-generated edit script: updates Leve profit-calculation columns
# (H: currentAveragePrice, I: currentAveragePriceNQ, J: currentAveragePriceHQ,
#  K: LevePriceNQ, L: LevePriceHQ, M: LeveProfitNQ, N: LeveProfitHQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 120.125
$ws.Range("I4").Value = 132.66667
$ws.Range("J4").Value = 82.5
$ws.Range("K4").Value = 132.66667
$ws.Range("L4").Value = 82.5
$ws.Range("M4").Value = -18.66667000000001
$ws.Range("N4").Value = -310.5
$ws.Range("H9").Value = 165.83333
$ws.Range("I9").Value = 99
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 99
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = -838
$ws.Range("H40").Value = 1674.5454
$ws.Range("I40").Value = 1668
$ws.Range("J40").Value = 1676.4706
$ws.Range("K40").Value = 1668
$ws.Range("L40").Value = 1676.4706
$ws.Range("M40").Value = -1493
$ws.Range("N40").Value = -2026.4706
$ws.Range("H86").Value = 78854.08
$ws.Range("I86").Value = 112867
$ws.Range("J86").Value = 2325
$ws.Range("K86").Value = 112867
$ws.Range("L86").Value = 2325
$ws.Range("M86").Value = -111744
$ws.Range("N86").Value = -4571
$ws.Range("H89").Value = 78854.08
$ws.Range("I89").Value = 112867
$ws.Range("J89").Value = 2325
$ws.Range("K89").Value = 564335
$ws.Range("L89").Value = 11625
$ws.Range("M89").Value = -558719
$ws.Range("N89").Value = -22857
$ws.Range("H92").Value = 1296.9474
$ws.Range("I92").Value = 1594.7
$ws.Range("J92").Value = 966.1111
$ws.Range("K92").Value = 1594.7
$ws.Range("L92").Value = 966.1111
$ws.Range("M92").Value = -346.7
$ws.Range("N92").Value = -3462.1111
$ws.Range("H113").Value = 3410
$ws.Range("I113").Value = 3133.3333
$ws.Range("J113").Value = 3825
$ws.Range("K113").Value = 3133.3333
$ws.Range("L113").Value = 3825
$ws.Range("M113").Value = 120.6667000000002
$ws.Range("N113").Value = -10333
$ws.Range("H129").Value = 27215.158
$ws.Range("J129").Value = 41082.44
$ws.Range("L129").Value = 123247.32
$ws.Range("N129").Value = -133247.32
$ws.Range("H137").Value = 1191.9565
$ws.Range("I137").Value = 1007.625
$ws.Range("K137").Value = 3022.875
$ws.Range("M137").Value = -472.875
$ws.Range("H138").Value = 1642.4728
$ws.Range("I138").Value = 1258.138
$ws.Range("J138").Value = 2071.1538
$ws.Range("K138").Value = 3774.414
$ws.Range("L138").Value = 6213.4614
$ws.Range("M138").Value = 1365.586
$ws.Range("N138").Value = -16493.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 333.33334
$ws.Range("J4").Value = 400
$ws.Range("L4").Value = 400
$ws.Range("N4").Value = -632
$ws.Range("H54").Value = 9780
$ws.Range("J54").Value = 9780
$ws.Range("L54").Value = 9780
$ws.Range("N54").Value = -11318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 17831.455
$ws.Range("I102").Value = 9794
$ws.Range("J102").Value = 54000
$ws.Range("K102").Value = 9794
$ws.Range("L102").Value = 54000
$ws.Range("M102").Value = -6549
$ws.Range("N102").Value = -60490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8.5
$ws.Range("I7").Value = 7.571429
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 7.571429
$ws.Range("L7").Value = 15
$ws.Range("M7").Value = 105.428571
$ws.Range("N7").Value = -241
$ws.Range("H31").Value = 5052719
$ws.Range("I31").Value = 2313.4375
$ws.Range("J31").Value = 9806041
$ws.Range("K31").Value = 2313.4375
$ws.Range("L31").Value = 9806041
$ws.Range("M31").Value = -2018.4375
$ws.Range("N31").Value = -9806631
$ws.Range("H34").Value = 5052719
$ws.Range("I34").Value = 2313.4375
$ws.Range("J34").Value = 9806041
$ws.Range("K34").Value = 2313.4375
$ws.Range("L34").Value = 9806041
$ws.Range("M34").Value = -2111.4375
$ws.Range("N34").Value = -9806445
$ws.Range("H132").Value = 2143.7437
$ws.Range("I132").Value = 1537.8125
$ws.Range("J132").Value = 4913.7144
$ws.Range("K132").Value = 4613.4375
$ws.Range("L132").Value = 14741.1432
$ws.Range("M132").Value = -2083.4375
$ws.Range("N132").Value = -19801.1432
$ws.Range("H134").Value = 1336.1428
$ws.Range("I134").Value = 1289.6923
$ws.Range("J134").Value = 1411.625
$ws.Range("K134").Value = 3869.0769
$ws.Range("L134").Value = 4234.875
$ws.Range("M134").Value = -1334.0769
$ws.Range("N134").Value = -9304.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 8931.9
$ws.Range("I33").Value = 11151.125
$ws.Range("J33").Value = 55
$ws.Range("K33").Value = 66906.75
$ws.Range("L33").Value = 330
$ws.Range("M33").Value = -66623.75
$ws.Range("N33").Value = -896
$ws.Range("H51").Value = 462.5
$ws.Range("I51").Value = 462.5
$ws.Range("K51").Value = 1387.5
$ws.Range("M51").Value = -927.5
$ws.Range("H99").Value = 10163.75
$ws.Range("I99").Value = 993.2
$ws.Range("K99").Value = 2979.6
$ws.Range("M99").Value = -733.6000000000004
$ws.Range("H131").Value = 3178546
$ws.Range("I131").Value = 18675
$ws.Range("J131").Value = 3832312.5
$ws.Range("K131").Value = 56025
$ws.Range("L131").Value = 11496937.5
$ws.Range("M131").Value = -50985
$ws.Range("N131").Value = -11507017.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2349.0625
$ws.Range("I126").Value = 2914.1428
$ws.Range("J126").Value = 1909.5555
$ws.Range("K126").Value = 8742.4284
$ws.Range("L126").Value = 5728.666499999999
$ws.Range("M126").Value = -6272.428400000001
$ws.Range("N126").Value = -10668.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1422.4445
$ws.Range("I22").Value = 2559.6
$ws.Range("J22").Value = 985.0769
$ws.Range("K22").Value = 2559.6
$ws.Range("L22").Value = 985.0769
$ws.Range("M22").Value = -2264.6
$ws.Range("N22").Value = -1575.0769
$ws.Range("H27").Value = 1422.4445
$ws.Range("I27").Value = 2559.6
$ws.Range("J27").Value = 985.0769
$ws.Range("K27").Value = 2559.6
$ws.Range("L27").Value = 985.0769
$ws.Range("M27").Value = -2452.6
$ws.Range("N27").Value = -1199.0769
$ws.Range("H46").Value = 1125.7826
$ws.Range("I46").Value = 387.625
$ws.Range("J46").Value = 1519.4667
$ws.Range("K46").Value = 387.625
$ws.Range("L46").Value = 1519.4667
$ws.Range("M46").Value = -199.625
$ws.Range("N46").Value = -1895.4667
$ws.Range("H133").Value = 22167.666
$ws.Range("J133").Value = 22167.666
$ws.Range("L133").Value = 22167.666
$ws.Range("N133").Value = -27227.666
